# Adds a new "2022" column (S) to the Financial Soundness Indicators table,
# restyles the existing "2021" column (R) to match the rest of the year
# columns, fills in the new values for 2022, and moves the active selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlPasteFormats = -4122   # xlPasteFormats

# --- Restyle column R (2021) so it matches the other year columns (Q, etc.) ---
$ws.Range("Q4").Copy()
$ws.Range("R4").PasteSpecial($xlPasteFormats)

$ws.Range("Q5").Copy()
$ws.Range("R5").PasteSpecial($xlPasteFormats)

$ws.Range("Q6").Copy()
$ws.Range("R6").PasteSpecial($xlPasteFormats)

$ws.Range("Q7").Copy()
$ws.Range("R7").PasteSpecial($xlPasteFormats)

$excel.CutCopyMode = $false

# --- Add the new 2022 column (S), copying formatting from column R ---
$ws.Range("R4").Copy()
$ws.Range("S4").PasteSpecial($xlPasteFormats)

$ws.Range("R5").Copy()
$ws.Range("S5").PasteSpecial($xlPasteFormats)

$ws.Range("R6").Copy()
$ws.Range("S6").PasteSpecial($xlPasteFormats)

$ws.Range("R7").Copy()
$ws.Range("S7").PasteSpecial($xlPasteFormats)

$excel.CutCopyMode = $false

# --- Values for the new 2022 column ---
$ws.Range("S4").Value = 2022
$ws.Range("S5").Value = 49.7
$ws.Range("S6").Value = 34.9
$ws.Range("S7").Value = 21

# --- Move the active selection like in the authored workbook ---
$ws.Range("R12").Select() | Out-Null
